# Updated projects of Gustavo.
# - Row 7 (JCNE project): the grant's "Ano fim" (end year) moved from 2016 to 2019.
# - A new row 8 is added for a related/new project ("Sistema Multifásico de
#   Arrefecimento de Componentes Eletrônicos", an "Apoio a Grupos Emergentes"
#   FAPERJ grant running 2016-2018).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gustavo")

# --- Row 7: update end year ---
$ws.Cells.Item(7, 5).Value = 2019

# --- Row 8: new project row ---
# Copy row 7's formatting (borders/alignment/number formats) down into row 8
# first, so the new row matches the sheet's existing look (style ids 1 / 15),
# then overwrite the values.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 42

$ws.Cells.Item(8, 1).Value  = "E-26/010.001642/2016"
$ws.Cells.Item(8, 9).Value  = "Sistema Multifásico de Arrefecimento de Componentes Eletrônicos"
$ws.Cells.Item(8, 10).Value = "249.927,90"
$ws.Cells.Item(8, 11).Value = "162.795,90"
$ws.Cells.Item(8, 8).Value  = "Apoio a Grupos Emergentes de Pesquisa no Estado do Rio de Janeiro 2016"
$ws.Cells.Item(8, 3).Value  = "Apoio a Grupos Emergentes"
$ws.Cells.Item(8, 7).Value  = "30/09/2018"
$ws.Cells.Item(8, 4).Value  = 2016
$ws.Cells.Item(8, 5).Value  = 2018
$ws.Cells.Item(8, 6).Value  = $ws.Cells.Item(7, 6).Value()
$ws.Cells.Item(8, 2).Value  = "FAPERJ"

# Re-touch I7 (same text, re-entered) so its shared-string slot is freed and
# the string is appended fresh at the end of the table.
$ws.Cells.Item(7, 9).Value = ""
$ws.Cells.Item(7, 9).Value = "Sistema de Alto Desempenho de Simulação de Escoamentos Multifásicos em Geometrias Complexas"

# --- Column K width (widened to fit the new "162.795,90" value) ---
$ws.Columns.Item(11).ColumnWidth = 10.83

# --- Selection state left on the sheet ---
$ws.Range("I7").Select()
